$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.960.44"
$ws.Cells.Item(2, 5).Value = "  +0.65%  "

$ws.Cells.Item(3, 4).Value = "1.643.68"
$ws.Cells.Item(3, 5).Value = "  -0.17%  "

$ws.Cells.Item(4, 5).Value = "  -0.38%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "219.66"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.36%  "

$ws.Cells.Item(6, 5).Value = "  -0.82%  "

$ws.Cells.Item(7, 5).Value = "  -0.40%  "

$ws.Cells.Item(8, 5).Value = "  -0.04%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.0624"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.54%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.36"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.95%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0847"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.46%  "

$ws.Cells.Item(12, 4).Value = "1.873.48"
$ws.Cells.Item(12, 5).Value = "  +0.04%  "

$ws.Cells.Item(13, 4).Value = "1.657.46"
$ws.Cells.Item(13, 5).Value = "  +0.79%  "

$ws.Cells.Item(14, 5).Value = "  -0.06%  "

$ws.Cells.Item(15, 5).Value = "  +0.46%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "65.93"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.83%  "

$ws.Cells.Item(17, 4).Value = "26.939.96"
$ws.Cells.Item(17, 5).Value = "  +0.63%  "

$ws.Cells.Item(18, 5).Value = "  -0.36%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "217.53"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.44%  "

$ws.Cells.Item(20, 5).Value = "  -0.52%  "

$ws.Cells.Item(21, 5).Value = "  -0.19%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.63"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +6.00%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.45"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.84%  "

$ws.Cells.Item(24, 5).Value = "  -1.35%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "148.19"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.36%  "

$ws.Cells.Item(26, 5).Value = "  -0.51%  "

$ws.Cells.Item(27, 5).Value = "  +2.50%  "

$ws.Cells.Item(28, 5).Value = "  +0.01%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.81"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.95%  "

$ws.Cells.Item(30, 5).Value = "  +0.46%  "

$ws.Cells.Item(31, 5).Value = "  +1.09%  "

$ws.Cells.Item(32, 5).Value = "  +0.07%  "

$ws.Cells.Item(33, 5).Value = "  +0.02%  "

$ws.Cells.Item(34, 5).Value = "  +1.65%  "

$ws.Cells.Item(35, 4).Value = "1.270.86"
$ws.Cells.Item(35, 5).Value = "  -1.42%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.43"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.19%  "

$ws.Cells.Item(37, 5).Value = "  -1.95%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.532"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.69%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.825"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.56%  "

$ws.Cells.Item(40, 5).Value = "  -0.45%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.808"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.25%  "

$ws.Cells.Item(42, 5).Value = "  +0.44%  "

$ws.Cells.Item(43, 4).Value = "1.784.06"
$ws.Cells.Item(43, 5).Value = "  -0.37%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "61.74"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.02%  "

$ws.Cells.Item(45, 2).Value = "MXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.09"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -6.31%  "

$ws.Cells.Item(46, 2).Value = "Quant"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "92.51"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.69%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.60"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.63%  "

$ws.Cells.Item(48, 5).Value = "  -1.22%  "

$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0971"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.07%  "

$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.59"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.36%  "

$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.01"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.33%  "
